$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row for "2021Q4" (old row 5) was removed; all following rows shift up by one.
$ws.Rows.Item(5).Delete()

# Re-write the data (values were recalculated after the row removal).
$ws.Cells.Item(2, 1).Value = "2020Q4"
$ws.Cells.Item(2, 2).Value = 0.108108108
$ws.Cells.Item(2, 3).Value = 0.4
$ws.Cells.Item(2, 4).Value = 0.008333333300000001
$ws.Cells.Item(2, 5).Value = -0.291891892
$ws.Cells.Item(2, 6).ClearContents()

$ws.Cells.Item(3, 1).Value = "2021Q1"
$ws.Cells.Item(3, 2).Value = 0.05704114339296366
$ws.Cells.Item(3, 3).Value = 0.3558812823004592
$ws.Cells.Item(3, 4).Value = 0.0749999999
$ws.Cells.Item(3, 5).Value = -0.2988401389074955
$ws.Cells.Item(3, 6).ClearContents()

$ws.Cells.Item(4, 1).Value = "2021Q2"
$ws.Cells.Item(4, 2).Value = 0.05111816859879303
$ws.Cells.Item(4, 3).Value = 0.3875183639680647
$ws.Cells.Item(4, 4).Value = 0.09166666650000001
$ws.Cells.Item(4, 5).Value = -0.3364001953692717
$ws.Cells.Item(4, 6).ClearContents()

$ws.Cells.Item(5, 1).Value = "2022Q1"
$ws.Cells.Item(5, 2).Value = 0.102908277
$ws.Cells.Item(5, 3).Value = 0.402684564
$ws.Cells.Item(5, 4).Value = 0.008333333300000001
$ws.Cells.Item(5, 5).Value = -0.299776287
$ws.Cells.Item(5, 6).ClearContents()

$ws.Cells.Item(6, 1).Value = "2022Q2"
$ws.Cells.Item(6, 2).Value = 0.04639970480210889
$ws.Cells.Item(6, 3).Value = 0.3888433962409018
$ws.Cells.Item(6, 4).Value = 0.1333333332
$ws.Cells.Item(6, 5).Value = -0.3424436914387929
$ws.Cells.Item(6, 6).Value = 0.1731867202354249

$ws.Cells.Item(7, 1).Value = "2022Q3"
$ws.Cells.Item(7, 2).Value = 0.05213945032610409
$ws.Cells.Item(7, 3).Value = 0.3773580867036184
$ws.Cells.Item(7, 4).Value = 0.2083333331
$ws.Cells.Item(7, 5).Value = -0.3252186363775144
$ws.Cells.Item(7, 6).Value = 0.08826959312244265

$ws.Cells.Item(8, 1).Value = "2022Q4"
$ws.Cells.Item(8, 2).Value = 0.07456890343205831
$ws.Cells.Item(8, 3).Value = 0.369937611989041
$ws.Cells.Item(8, 4).Value = 1.8884928218
$ws.Cells.Item(8, 5).Value = -0.2953687085569827
$ws.Cells.Item(8, 6).Value = -0.1219722443004176

$ws.Cells.Item(9, 1).Value = "2023Q1"
$ws.Cells.Item(9, 2).Value = 0.07838743747622531
$ws.Cells.Item(9, 3).Value = 0.3639284472021985
$ws.Cells.Item(9, 4).Value = 1.4791228062
$ws.Cells.Item(9, 5).Value = -0.2855410097259732
$ws.Cells.Item(9, 6).Value = -0.04748633528183921

$ws.Cells.Item(10, 1).Value = "2023Q2"
$ws.Cells.Item(10, 2).Value = 0.07957850659290683
$ws.Cells.Item(10, 3).Value = 0.3703992569229075
$ws.Cells.Item(10, 4).Value = 2.4665909069
$ws.Cells.Item(10, 5).Value = -0.2908207503300007
$ws.Cells.Item(10, 6).Value = -0.1507486994194464

$ws.Cells.Item(11, 1).Value = "2023Q3"
$ws.Cells.Item(11, 2).Value = 0.07757506325723942
$ws.Cells.Item(11, 3).Value = 0.3599759990335085
$ws.Cells.Item(11, 4).Value = 2.5035964888
$ws.Cells.Item(11, 5).Value = -0.282400935776269
$ws.Cells.Item(11, 6).Value = -0.1316582010126335

$ws.Cells.Item(12, 1).Value = "2023Q4"
$ws.Cells.Item(12, 2).Value = 0.07626219378725912
$ws.Cells.Item(12, 3).Value = 0.3634170471661184
$ws.Cells.Item(12, 4).Value = 2.7082456115
$ws.Cells.Item(12, 5).Value = -0.2871548533788593
$ws.Cells.Item(12, 6).Value = -0.02780881975701488

$ws.Cells.Item(13, 1).Value = "2024Q1"
$ws.Cells.Item(13, 2).Value = 0.07756450574887219
$ws.Cells.Item(13, 3).Value = 0.3605270496924279
$ws.Cells.Item(13, 4).Value = 2.7858492797
$ws.Cells.Item(13, 5).Value = -0.2829625439435557
$ws.Cells.Item(13, 6).Value = -0.009030106690775974

$ws.Cells.Item(14, 1).Value = "2024Q2"
$ws.Cells.Item(14, 2).Value = 0.07650396207711023
$ws.Cells.Item(14, 3).Value = 0.3490339542141579
$ws.Cells.Item(14, 4).Value = 2.6880143514
$ws.Cells.Item(14, 5).Value = -0.2725299921370476
$ws.Cells.Item(14, 6).Value = -0.06289358022836444

$ws.Cells.Item(15, 1).Value = "2024Q3"
$ws.Cells.Item(15, 2).Value = 0.07587408062511565
$ws.Cells.Item(15, 3).Value = 0.343827116858486
$ws.Cells.Item(15, 4).Value = 2.6458732031
$ws.Cells.Item(15, 5).Value = -0.2679530362333704
$ws.Cells.Item(15, 6).Value = -0.05116094783179104

$ws.Cells.Item(16, 1).Value = "2024Q4"
$ws.Cells.Item(16, 2).Value = 0.07856301634196641
$ws.Cells.Item(16, 3).Value = 0.3458978321782535
$ws.Cells.Item(16, 4).Value = 2.6684808587
$ws.Cells.Item(16, 5).Value = -0.267334815836287
$ws.Cells.Item(16, 6).Value = -0.06902212276531694

$ws.Cells.Item(17, 1).Value = "2025Q1"
$ws.Cells.Item(17, 2).Value = 0.07926348720140522
$ws.Cells.Item(17, 3).Value = 0.3414053954863104
$ws.Cells.Item(17, 4).Value = 2.3345295035
$ws.Cells.Item(17, 5).Value = -0.2621419082849052
$ws.Cells.Item(17, 6).Value = -0.07358088942967556

$ws.Cells.Item(18, 1).Value = "2025Q2"
$ws.Cells.Item(18, 2).Value = 0.08132440800948877
$ws.Cells.Item(18, 3).Value = 0.3400601736860327
$ws.Cells.Item(18, 4).Value = 2.7126315763
$ws.Cells.Item(18, 5).Value = -0.2587357656765439
$ws.Cells.Item(18, 6).Value = -0.05061544365204029

$ws.Cells.Item(19, 1).Value = "2025Q3"
$ws.Cells.Item(19, 2).Value = 0.08223145772514558
$ws.Cells.Item(19, 3).Value = 0.3365925883083096
$ws.Cells.Item(19, 4).Value = 2.4500159467
$ws.Cells.Item(19, 5).Value = -0.254361130583164
$ws.Cells.Item(19, 6).Value = -0.05072495479531958

$ws.Cells.Item(20, 1).Value = "2025Q4"
$ws.Cells.Item(20, 2).Value = 0.07916596409077691
$ws.Cells.Item(20, 3).Value = 0.3306924009905033
$ws.Cells.Item(20, 4).Value = 1.8167703333
$ws.Cells.Item(20, 5).Value = -0.2515264368997264
$ws.Cells.Item(20, 6).Value = -0.05913325911968592

$ws.Cells.Item(21, 1).Value = "2026Q1"
$ws.Cells.Item(21, 2).Value = 0.07964367563289203
$ws.Cells.Item(21, 3).Value = 0.3288143450685608
$ws.Cells.Item(21, 4).Value = 1.1324202542
$ws.Cells.Item(21, 5).Value = -0.2491706694356687
$ws.Cells.Item(21, 6).Value = -0.04948174419764617
